$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.323.76'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '3.282.87'
$ws.Range('E3').Value = '  -2.52%  '
$c = $ws.Range('D4')
$c.Value = "'0.998"
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '
$c = $ws.Range('D5')
$c.Value = "'580.52"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.42%  '
$c = $ws.Range('D6')
$c.Value = "'174.75"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -7.28%  '
$c = $ws.Range('D7')
$c.Value = "'0.996"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').Value = '3.280.96'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('E10').Value = '  -5.37%  '
$c = $ws.Range('D11')
$c.Value = "'0.572"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -2.52%  '
$c = $ws.Range('D12')
$c.Value = "'45.14"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -5.33%  '
$c = $ws.Range('D13')
$c.Value = "'0.0000271"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.30%  '
$c = $ws.Range('D14')
$c.Value = "'667.15"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '3.805.17'
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('E16').Value = '  -3.99%  '
$ws.Range('D17').Value = '67.279.96'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').Value = '3.270.01'
$ws.Range('E19').Value = '  -3.02%  '
$c = $ws.Range('D20')
$c.Value = "'17.33"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -4.09%  '
$c = $ws.Range('D21')
$c.Value = "'10.78"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -3.67%  '
$c = $ws.Range('D22')
$c.Value = "'0.881"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D23')
$c.Value = "'17.00"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -5.87%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D24')
$c.Value = "'5.33"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +4.38%  '
$c = $ws.Range('D25')
$c.Value = "'97.51"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -3.34%  '
$ws.Range('E26').Value = '  -4.01%  '
$c = $ws.Range('D27')
$c.Value = "'2.66"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -6.58%  '
$c = $ws.Range('D28')
$c.Value = "'9.22"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -5.72%  '
$c = $ws.Range('D29')
$c.Value = "'32.57"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.85%  '
$c = $ws.Range('D30')
$c.Value = "'8.33"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -4.66%  '
$c = $ws.Range('D31')
$c.Value = "'6.92"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.01%  '
$c = $ws.Range('D32')
$c.Value = "'568.13"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -8.17%  '
$ws.Range('E33').Value = '  -2.90%  '
$ws.Range('D34').Value = '3.758.73'
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('E35').Value = '  -3.91%  '
$c = $ws.Range('D36')
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.20%  '
$c = $ws.Range('D37')
$c.Value = "'3.38"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -13.07%  '
$c = $ws.Range('D38')
$c.Value = "'55.58"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('E40').Value = '  -7.35%  '
$c = $ws.Range('D41')
$c.Value = "'32.29"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -4.34%  '
$c = $ws.Range('D42')
$c.Value = "'3.05"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -6.79%  '
$ws.Range('D43').Value = '0.0₃0662'
$ws.Range('E43').Value = '  -6.24%  '
$c = $ws.Range('D44')
$c.Value = "'0.327"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -5.52%  '
$c = $ws.Range('D45')
$c.Value = "'3.22"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.93%  '
$c = $ws.Range('D46')
$c.Value = "'0.0401"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -5.33%  '
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D48')
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D49')
$c.Value = "'0.126"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.43%  '
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('E51').Value = '  -4.98%  '
